$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 1281.3636
$ws.Range("I55").Value = 675.25
$ws.Range("K55").Value = 675.25
$ws.Range("M55").Value = -461.25
$ws.Range("H106").Value = 2500
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 2500
$ws.Range("K106").Value = 0
$ws.Range("L106").ClearContents()
$ws.Range("M106").Value = 2500
$ws.Range("N106").Value = -3762
$ws.Range("H111").Value = 1546.5
$ws.Range("I111").Value = 608.1667
$ws.Range("J111").Value = 1781.0834
$ws.Range("K111").Value = 1824.5001
$ws.Range("L111").Value = 5343.2502
$ws.Range("M111").Value = 1242.4999
$ws.Range("N111").Value = -11477.2502
$ws.Range("H135").Value = 2090.7144
$ws.Range("I135").Value = 1712.8
$ws.Range("K135").Value = 15415.2
$ws.Range("M135").Value = -12880.2
$ws.Range("H137").Value = 2218.8
$ws.Range("I137").Value = 1399.8334
$ws.Range("J137").Value = 2569.7856
$ws.Range("K137").Value = 4199.5002
$ws.Range("L137").Value = 7709.3568
$ws.Range("M137").Value = -1649.5002
$ws.Range("N137").Value = -12809.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = 0
$ws.Range("H88").Value = 1307.3077
$ws.Range("J88").Value = 968.5
$ws.Range("L88").Value = 968.5
$ws.Range("N88").Value = -1780.5
$ws.Range("H91").Value = 1307.3077
$ws.Range("J91").Value = 968.5
$ws.Range("L91").Value = 968.5
$ws.Range("N91").Value = -3776.5
$ws.Range("H102").Value = 13895867
$ws.Range("J102").Value = 9779.799999999999
$ws.Range("L102").Value = 9779.799999999999
$ws.Range("N102").Value = -13023.8
$ws.Range("H110").Value = 125002630
$ws.Range("I110").Value = 200001020
$ws.Range("J110").Value = 5314.3335
$ws.Range("K110").Value = 200001020
$ws.Range("L110").Value = 5314.3335
$ws.Range("M110").Value = -199998975
$ws.Range("N110").Value = -9404.333500000001
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").ClearContents()
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = 0
$ws.Range("H122").Value = 948.875
$ws.Range("I122").Value = 824.4167
$ws.Range("J122").Value = 1322.25
$ws.Range("K122").Value = 2473.2501
$ws.Range("L122").Value = 3966.75
$ws.Range("M122").Value = -23.2501000000002
$ws.Range("N122").Value = -8866.75
$ws.Range("H132").Value = 1043.2727
$ws.Range("I132").Value = 1116.9445
$ws.Range("K132").Value = 3350.8335
$ws.Range("M132").Value = -820.8335000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = 0
$ws.Range("H20").Value = 4999
$ws.Range("J20").Value = 4999
$ws.Range("L20").Value = 4999
$ws.Range("N20").Value = -5493

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4595.387
$ws.Range("I31").Value = 1454.9166
$ws.Range("K31").Value = 1454.9166
$ws.Range("M31").Value = -1159.9166
$ws.Range("H34").Value = 4595.387
$ws.Range("I34").Value = 1454.9166
$ws.Range("K34").Value = 1454.9166
$ws.Range("M34").Value = -1252.9166
$ws.Range("H99").Value = 2043.9231
$ws.Range("I99").Value = 2211.8572
$ws.Range("J99").Value = 1848
$ws.Range("K99").Value = 2211.8572
$ws.Range("L99").Value = 1848
$ws.Range("M99").Value = -713.8571999999999
$ws.Range("N99").Value = -4844
$ws.Range("H105").Value = 2864.9
$ws.Range("I105").Value = 1347.6
$ws.Range("J105").Value = 4382.2
$ws.Range("K105").Value = 1347.6
$ws.Range("L105").Value = 4382.2
$ws.Range("M105").Value = 399.4000000000001
$ws.Range("N105").Value = -7876.2
$ws.Range("H126").Value = 2043.9231
$ws.Range("I126").Value = 2211.8572
$ws.Range("J126").Value = 1848
$ws.Range("K126").Value = 6635.571599999999
$ws.Range("L126").Value = 5544
$ws.Range("M126").Value = -4165.571599999999
$ws.Range("N126").Value = -10484
$ws.Range("H132").Value = 2114.6365
$ws.Range("I132").Value = 2001.1
$ws.Range("K132").Value = 6003.299999999999
$ws.Range("M132").Value = -3473.299999999999
$ws.Range("H134").Value = 3981.3333
$ws.Range("I134").Value = 3981.3333
$ws.Range("K134").Value = 11943.9999
$ws.Range("M134").Value = -9408.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1763.7142
$ws.Range("I5").Value = 1198.8334
$ws.Range("K5").Value = 3596.5002
$ws.Range("M5").Value = -3484.5002
$ws.Range("H42").Value = 3000
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H132").Value = 3613.4614
$ws.Range("I132").Value = 3388.7144
$ws.Range("J132").Value = 3875.6667
$ws.Range("K132").Value = 30498.4296
$ws.Range("L132").Value = 34881.0003
$ws.Range("M132").Value = -27968.4296
$ws.Range("N132").Value = -39941.0003
$ws.Range("H135").Value = 1763.7142
$ws.Range("I135").Value = 1198.8334
$ws.Range("K135").Value = 10789.5006
$ws.Range("M135").Value = -8254.500599999999
$ws.Range("H137").Value = 4912.857
$ws.Range("I137").Value = 2800
$ws.Range("K137").Value = 8400
$ws.Range("M137").Value = -3300

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 180.75
$ws.Range("I9").Value = 200.14285
$ws.Range("J9").Value = 45
$ws.Range("K9").Value = 200.14285
$ws.Range("L9").Value = 45
$ws.Range("M9").Value = -30.14285000000001
$ws.Range("N9").Value = -385
$ws.Range("H97").Value = 758.8333
$ws.Range("I97").Value = 339.75
$ws.Range("J97").Value = 968.375
$ws.Range("K97").Value = 339.75
$ws.Range("L97").Value = 968.375
$ws.Range("M97").Value = 156.25
$ws.Range("N97").Value = -1960.375
$ws.Range("H122").Value = 1968.3125
$ws.Range("I122").Value = 1373
$ws.Range("K122").Value = 4119
$ws.Range("M122").Value = -1669
$ws.Range("H132").Value = 1389.8846
$ws.Range("I132").Value = 1421.68
$ws.Range("K132").Value = 4265.04
$ws.Range("M132").Value = -1735.04

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = 0
$ws.Range("H17").Value = 6000
$ws.Range("I17").Value = 6000
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -5830
$ws.Range("H22").Value = 525
$ws.Range("I22").Value = 362.375
$ws.Range("J22").Value = 850.25
$ws.Range("K22").Value = 362.375
$ws.Range("L22").Value = 850.25
$ws.Range("M22").Value = -67.375
$ws.Range("N22").Value = -1440.25
$ws.Range("H27").Value = 525
$ws.Range("I27").Value = 362.375
$ws.Range("J27").Value = 850.25
$ws.Range("K27").Value = 362.375
$ws.Range("L27").Value = 850.25
$ws.Range("M27").Value = -255.375
$ws.Range("N27").Value = -1064.25
$ws.Range("H40").Value = 2640.5454
$ws.Range("I40").Value = 3413.1428
$ws.Range("K40").Value = 3413.1428
$ws.Range("M40").Value = -3277.1428
$ws.Range("H55").Value = 1340.5714
$ws.Range("I55").Value = 1420.8334
$ws.Range("J55").Value = 1280.375
$ws.Range("K55").Value = 1420.8334
$ws.Range("L55").Value = 1280.375
$ws.Range("M55").Value = -1247.8334
$ws.Range("N55").Value = -1626.375
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("N125").Value = 0
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = 0
$ws.Range("H132").Value = 2532.6428
$ws.Range("I132").Value = 2450.6365
$ws.Range("J132").Value = 2833.3333
$ws.Range("K132").Value = 7351.9095
$ws.Range("L132").Value = 8499.999899999999
$ws.Range("M132").Value = -4821.9095
$ws.Range("N132").Value = -13559.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 1675.5
$ws.Range("I3").Value = 2475
$ws.Range("K3").Value = 2475
$ws.Range("M3").Value = -2361
$ws.Range("H81").Value = 685.8
$ws.Range("I81").Value = 685.8
$ws.Range("K81").Value = 1371.6
$ws.Range("M81").Value = -310.5999999999999
$ws.Range("H84").Value = 685.8
$ws.Range("I84").Value = 685.8
$ws.Range("K84").Value = 6858
$ws.Range("M84").Value = -1554
$ws.Range("H122").Value = 2856.4375
$ws.Range("I122").Value = 2661.8462
$ws.Range("J122").Value = 3699.6667
$ws.Range("K122").Value = 7985.5386
$ws.Range("L122").Value = 11099.0001
$ws.Range("M122").Value = -5535.5386
$ws.Range("N122").Value = -15999.0001
$ws.Range("H126").Value = 4920.5
$ws.Range("I126").Value = 1815.5
$ws.Range("J126").Value = 7249.25
$ws.Range("K126").Value = 5446.5
$ws.Range("L126").Value = 21747.75
$ws.Range("M126").Value = -2976.5
$ws.Range("N126").Value = -26687.75
$ws.Range("H132").Value = 998
$ws.Range("J132").Value = 1666.6666
$ws.Range("L132").Value = 4999.9998
$ws.Range("N132").Value = -10059.9998
$ws.Range("H136").Value = 2861.4666
$ws.Range("I136").Value = 1998.3
$ws.Range("K136").Value = 5994.9
$ws.Range("M136").Value = -3444.9
